$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values for columns B-E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update CON row (row 2) values for columns B-E
$ws.Range("B2").Value = -1.9546825443175493
$ws.Range("C2").Value = 4.397415285381463
$ws.Range("D2").Value = 3.3322129417138484
$ws.Range("E2").Value = 6.4394897193417506

# Update STR row (row 3) values for columns B-E
$ws.Range("B3").Value = 3.6266888763321674
$ws.Range("C3").Value = 12.123617845552914
$ws.Range("D3").Value = 15.967788290335697
$ws.Range("E3").Value = 4.0304374520007684

# Update the selected range to match the narrower selection
$ws.Range("B1:E3").Select()
